# "Updated symbol list" refresh (GitHub Actions scrape) — Fri Dec 30 03:53:33 UTC 2022
#
# Updates a batch of Price (col D) values, plus a name/link/rank reshuffle
# for three coins that swapped rank positions (rows 41-43), and two
# "Best/Worst in 24h" label moves (E18, E41-E43).
#
# Column D cells store plain numeric-looking text (e.g. "245.80", "0.001514")
# as literal strings, not numbers (no number formatting is applied in the
# sheet). Assigning a bare numeric-looking string to Range.Value lets Excel's
# normal "smart" entry coerce it into a real number (and drop formatting
# such as trailing zeros), so every column-D write below is prefixed with a
# leading apostrophe, exactly as a user would type it in the Excel UI to
# force literal text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$q  = "'"   # leading-apostrophe text-prefix, prepended to numeric-looking strings

# Row 2 - BNB
$ws.Range("D2").Value = $q + "245.80"

# Row 3 - OKB
$ws.Range("D3").Value = $q + "24.11"

# Row 4 - HuobiToken
$ws.Range("D4").Value = $q + "5.253"

# Row 5 - Cronos
$ws.Range("D5").Value = $q + "0.05785"

# Row 6 - KuCoinToken
$ws.Range("D6").Value = $q + "6.500"

# Row 7 - GateToken
$ws.Range("D7").Value = $q + "3.142"

# Row 9 - FTXToken
$ws.Range("D9").Value = $q + "0.8565"

# Row 10 - WazirX
$ws.Range("D10").Value = $q + "0.1364"

# Row 11 - MandalaExchangeToken
$ws.Range("D11").Value = $q + "0.06945"

# Row 12 - LiechtensteinCryptoassetsExchange
$ws.Range("D12").Value = $q + "0.03186"

# Row 13 - BitrueCoin
$ws.Range("D13").Value = $q + "0.02878"

# Row 14 - BitMartToken
$ws.Range("D14").Value = $q + "0.09391"

# Row 15 - MCDex
$ws.Range("D15").Value = $q + "3.752"

# Row 16 - BitForexToken
$ws.Range("D16").Value = $q + "0.001514"

# Row 18 - One (price refresh + gains the "Worst in 24h" tag)
$ws.Range("D18").Value = $q + "0.0005992"
$ws.Range("E18").Value = "17OneONEWorstin24h"

# Row 19 - TigerCash
$ws.Range("D19").Value = $q + "0.006270"

# Row 20 - BitKan
$ws.Range("D20").Value = $q + "0.001237"

# Row 21 - HotbitToken
$ws.Range("D21").Value = $q + "0.004613"

# Row 22 - NitroEx
$ws.Range("D22").Value = $q + "0.00006901"

# Row 23 - LEO
$ws.Range("D23").Value = $q + "3.513"

# Row 24 - BTSEToken
$ws.Range("D24").Value = $q + "2.116"

# Row 25 - BitpandaEcosystemToken
$ws.Range("D25").Value = $q + "0.3192"

# Row 26 - ProBitToken
$ws.Range("D26").Value = $q + "0.1347"

# Row 28 - UpBots
$ws.Range("D28").Value = $q + "0.0002330"

# Row 40 - IDEX
$ws.Range("D40").Value = $q + "0.03654"

# Rows 41-43 - BKEXToken / CEJI / KickToken cycle ranks: row 41 becomes
# KickToken, row 42 becomes BKEXToken, row 43 becomes CEJI (each keeping its
# own fresh price/link/rank-label for the new row).
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = $q + "0.006291"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = $q + "0.1054"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = $q + "0.003401"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"

# Row 44 - LocalTraders
$ws.Range("D44").Value = $q + "0.008059"

# Row 45 - CoinLion
$ws.Range("D45").Value = $q + "0.00005273"

# Row 47 - CoinbaseStockToken
$ws.Range("D47").Value = $q + "0.3501"

# Row 48 - BOLO
$ws.Range("D48").Value = $q + "0.002340"
